$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for the 380 kV case (Case_4_50).
# Maps each affected column letter to its 24 new values for rows 2-25.
$data = @{
    "B" = @(21.45275908512185, 20.62898603025496, 20.11508000833933, 19.90398743118845, 19.86884607542741, 20.11223939943191, 21.17061401664428, 23.16739397162576, 24.5696186092189, 25.1905517926124, 25.42304490274724, 25.37309375871779, 25.2097331371156, 25.10932054490632, 24.52868009389082, 24.16798164485526, 23.95893275441659, 23.88788689739929, 24.20654418539318, 25.25778927173745, 25.92937038129202, 25.57241064009423, 24.18911527334118, 22.63750400909523)
    "C" = @(20.02945465575656, 19.04399962648701, 18.42049241153628, 18.16215328699376, 18.11901152454517, 18.41702504311991, 19.69370680510226, 22.03702232003041, 23.6456102532911, 24.35044384574589, 24.61331962926647, 24.5568860447829, 24.37215235222047, 24.25846869016014, 23.59899179683778, 23.18740679508483, 22.94814745136317, 22.86670972516779, 23.23148344776928, 24.42652358873323, 25.18401135299923, 24.78192469513845, 23.21156458399724, 21.42200435912712)
    "D" = @(5.9957883131722, 5.882050128650533, 5.813374750592495, 5.785718620852172, 5.78114734959065, 5.813000385891694, 5.95635337360175, 6.24500516489932, 6.459462604706394, 6.557091458098065, 6.594035339273757, 6.58608053121655, 6.560131577130544, 6.544232653929467, 6.453080621661051, 6.397153684497246, 6.364994752853241, 6.354108928515162, 6.403106563495577, 6.567754401487505, 6.675193879697368, 6.61787825892306, 6.400415282165503, 6.166341179174973)
    "E" = @(16.37683494653892, 15.44446074142737, 14.8495677531223, 14.60178108072154, 14.56032140596913, 14.84624733039862, 16.06015146876109, 18.3083110563658, 19.97433918212398, 20.69188093802045, 20.95784268837326, 20.90081832834722, 20.7138768634388, 20.5986216047381, 19.92663822692864, 19.50409776133959, 19.25726367554778, 19.17303586137763, 19.54947087004511, 20.76894194900202, 21.53240917786331, 21.12798504964102, 19.52896986964194, 17.68119384690134)
    "G" = @(3.653486065520515, 3.659877664583931, 3.663986330843474, 3.665707249719177, 3.665995830090506, 3.664009350680308, 3.655651831531955, 3.640711167524708, 3.630598725155235, 3.626181981818574, 3.624535534189389, 3.624888971231239, 3.626046006709355, 3.62675811129149, 3.630891035471737, 3.633473223478136, 3.634975715993074, 3.635487412262185, 3.633196558283029, 3.625705452052286, 3.620961435393506, 3.623479613026344, 3.633321582659065, 3.644599855073086)
    "I" = @(28.54966882240423, 28.21839105510682, 28.02220408376729, 27.94413462771861, 27.93128644964328, 28.0211435266529, 28.43398519783494, 29.29786841194011, 29.96144686722331, 30.26864936202131, 30.38566472752774, 30.36043425888939, 30.27826304594565, 30.22801747056198, 29.9414707720103, 29.76699173510785, 29.6671434055733, 29.63342614574649, 29.78551341318455, 30.3023808105108, 30.64413446138937, 30.46140006888973, 29.77713831722494, 29.0587761718688)
    "N" = @(18.53878884512057, 18.56600056890706, 18.58465433105099, 18.5927414054744, 18.59411347355112, 18.58476143522236, 18.54776536611779, 18.4908241630176, 18.45874241492369, 18.4463178271214, 18.44192935317668, 18.44286034348776, 18.44595041681162, 18.44788452066793, 18.45959842819588, 18.46734336182201, 18.47200198038714, 18.473614200017, 18.46649776155139, 18.44503416304879, 18.43285309012255, 18.43918389464662, 18.46687941627363, 18.50453299374485)
}

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5; "G" = 7; "I" = 9; "N" = 14 }

foreach ($col in $data.Keys) {
    $colNum = $colIndex[$col]
    $values = $data[$col]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $rowNum = $i + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $values[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) columns across rows 2-25."
